$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tkibuli")

$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2023

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 1473.2

$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 645

$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 1688.7

$excel.CutCopyMode = $false
